$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 49
$ws.Cells.Item(2, 2).Value = 61
$ws.Cells.Item(2, 3).Value = 36
$ws.Cells.Item(2, 5).Value = 26
$ws.Cells.Item(2, 6).Value = 29
$ws.Cells.Item(2, 7).Value = 3
$ws.Cells.Item(2, 8).Value = 30
$ws.Cells.Item(2, 9).Value = 4
$ws.Cells.Item(3, 3).Value = 18
$ws.Cells.Item(4, 3).Value = 26
$ws.Cells.Item(5, 3).Value = 30
$ws.Cells.Item(6, 3).Value = 14
$ws.Cells.Item(7, 3).Value = 55
$ws.Cells.Item(8, 3).Value = 50
$ws.Cells.Item(9, 3).Value = 56
$ws.Cells.Item(10, 3).Value = 22
$ws.Cells.Item(11, 3).Value = 21
$ws.Cells.Item(12, 3).Value = 7
$ws.Cells.Item(13, 3).Value = 11
$ws.Cells.Item(14, 3).Value = 4
$ws.Cells.Item(15, 3).Value = 57
$ws.Cells.Item(16, 3).Value = 28
$ws.Cells.Item(17, 3).Value = 51
$ws.Cells.Item(18, 3).Value = 23
$ws.Cells.Item(19, 3).Value = 0
$ws.Cells.Item(20, 3).Value = 21
$ws.Cells.Item(21, 3).Value = 29
$ws.Cells.Item(22, 3).Value = 2
$ws.Cells.Item(23, 3).Value = 1
$ws.Cells.Item(24, 3).Value = 46
$ws.Cells.Item(25, 3).Value = 31
$ws.Cells.Item(27, 3).Value = 48
$ws.Cells.Item(28, 3).Value = 31
$ws.Cells.Item(29, 3).Value = 53
$ws.Cells.Item(30, 3).Value = 27
$ws.Cells.Item(31, 3).Value = 29
$ws.Cells.Item(32, 3).Value = 43
$ws.Cells.Item(33, 3).Value = 47
$ws.Cells.Item(34, 3).Value = 23
$ws.Cells.Item(35, 3).Value = 58
$ws.Cells.Item(36, 3).Value = 44
$ws.Cells.Item(37, 3).Value = 48
$ws.Cells.Item(38, 3).Value = 35
$ws.Cells.Item(39, 3).Value = 40
$ws.Cells.Item(40, 3).Value = 40
$ws.Cells.Item(41, 3).Value = 33
$ws.Cells.Item(42, 3).Value = 34
$ws.Cells.Item(43, 3).Value = 42
$ws.Cells.Item(44, 3).Value = 53
$ws.Cells.Item(45, 3).Value = 1
$ws.Cells.Item(46, 3).Value = 5
$ws.Cells.Item(47, 3).Value = 22
$ws.Cells.Item(48, 3).Value = 31
$ws.Cells.Item(49, 3).Value = 40
$ws.Cells.Item(50, 3).Value = 49

$ws.Range("A51:I91").EntireRow.Delete()
